$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report data between row 2 and row 3
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

# Save row 2 current values
$d2 = $ws.Cells.Item(2, 4).Value2   # D2 Fecha
$m2 = $ws.Cells.Item(2, 13).Value2  # M2 Volumen
$n2 = $ws.Cells.Item(2, 14).Value2  # N2 Precio minimo
$o2 = $ws.Cells.Item(2, 15).Value2  # O2 Precio maximo
$p2 = $ws.Cells.Item(2, 16).Value2  # P2 Precio promedio ponderado
$s2 = $ws.Cells.Item(2, 19).Value2  # S2 Precio $/Kg

# Save row 3 current values
$d3 = $ws.Cells.Item(3, 4).Value2
$m3 = $ws.Cells.Item(3, 13).Value2
$n3 = $ws.Cells.Item(3, 14).Value2
$o3 = $ws.Cells.Item(3, 15).Value2
$p3 = $ws.Cells.Item(3, 16).Value2
$s3 = $ws.Cells.Item(3, 19).Value2

# Write row 3's original values into row 2
$ws.Cells.Item(2, 4).Value2 = $d3
$ws.Cells.Item(2, 13).Value2 = $m3
$ws.Cells.Item(2, 14).Value2 = $n3
$ws.Cells.Item(2, 15).Value2 = $o3
$ws.Cells.Item(2, 16).Value2 = $p3
$ws.Cells.Item(2, 19).Value2 = $s3

# Write row 2's original values into row 3
$ws.Cells.Item(3, 4).Value2 = $d2
$ws.Cells.Item(3, 13).Value2 = $m2
$ws.Cells.Item(3, 14).Value2 = $n2
$ws.Cells.Item(3, 15).Value2 = $o2
$ws.Cells.Item(3, 16).Value2 = $p2
$ws.Cells.Item(3, 19).Value2 = $s2

$wb.Save()
